# Update "想去人数" (want-to-go count) figures in column F across the
# three sheets that carry this data: 展览(1), 演出(2) and the combined
# 全部类型(4) sheet, matching the refreshed scrape output.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 509
$ws.Range("F4").Value = 451
$ws.Range("F5").Value = 8734
$ws.Range("F6").Value = 15
$ws.Range("F7").Value = 11110
$ws.Range("F10").Value = 21
$ws.Range("F13").Value = 122
$ws.Range("F15").Value = 297
$ws.Range("F18").Value = 85
$ws.Range("F20").Value = 419
$ws.Range("F22").Value = 1891
$ws.Range("F23").Value = 705
$ws.Range("F24").Value = 627
$ws.Range("F25").Value = 357
$ws.Range("F26").Value = 292
$ws.Range("F28").Value = 603
$ws.Range("F30").Value = 1279
$ws.Range("F31").Value = 25
$ws.Range("F33").Value = 6
$ws.Range("F36").Value = 464
$ws.Range("F37").Value = 3
$ws.Range("F39").Value = 301
$ws.Range("F41").Value = 141
$ws.Range("F42").Value = 534
$ws.Range("F43").Value = 378
$ws.Range("F44").Value = 116
$ws.Range("F45").Value = 813
$ws.Range("F46").Value = 656
$ws.Range("F48").Value = 150
$ws.Range("F49").Value = 137

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F8").Value = 54
$ws.Range("F14").Value = 29
$ws.Range("F19").Value = 101

# Sheet 4: 全部类型 (combined view, mirrors the same rows updated above)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 509
$ws.Range("F7").Value = 451
$ws.Range("F8").Value = 8734
$ws.Range("F9").Value = 15
$ws.Range("F10").Value = 11110
$ws.Range("F12").Value = 21
$ws.Range("F14").Value = 122
$ws.Range("F15").Value = 297
$ws.Range("F17").Value = 85
$ws.Range("F19").Value = 1891
$ws.Range("F20").Value = 705
$ws.Range("F21").Value = 627
$ws.Range("F22").Value = 357
$ws.Range("F23").Value = 292
$ws.Range("F25").Value = 603
$ws.Range("F26").Value = 54
$ws.Range("F29").Value = 1279
$ws.Range("F30").Value = 25
$ws.Range("F32").Value = 6
$ws.Range("F33").Value = 29
$ws.Range("F38").Value = 464
$ws.Range("F41").Value = 534
$ws.Range("F42").Value = 378
$ws.Range("F43").Value = 116
$ws.Range("F46").Value = 656
$ws.Range("F48").Value = 150
$ws.Range("F49").Value = 137
